# Update the table style ("Medium Style 2 - Accent 1" -> another built-in
# table style) on the "Sources of finance" table found on slide 6.
#
# The table lives in the 2nd shape (a graphicFrame) of slide 6; it is the
# only table in the whole deck, so we locate it defensively by scanning the
# slide's shapes for HasTable rather than hard-coding the shape index.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.ApplyStyle("{3F0C99E5-47E1-474C-856D-C17D15043B8B}")
    }
}
